$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2,3,4,6,7,8: set B,C,D columns to 0
$ws.Range("B2:D4").Value = 0
$ws.Range("B6:D8").Value = 0

# Row 5: B5 -> 0, C5 and D5 -> new values
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 0.7376383316729367
$ws.Range("D5").Value = 0.6659414307212431

# Row 9: B9 -> 0, C9 and D9 -> new values
$ws.Range("B9").Value = 0
$ws.Range("C9").Value = 0.6926875232078059
$ws.Range("D9").Value = -0.7773164866342108
